# Update "想去人数" (number of people interested) figures that changed
# between the two data refreshes, on both the "展览" sheet and the
# "全部类型" sheet (which both list the same events, offset by two rows).

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-5, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1047
$wsExhibit.Range("F3").Value = 280
$wsExhibit.Range("F4").Value = 2728
$wsExhibit.Range("F5").Value = 54

# Sheet "全部类型": rows 4-7, column F (same events as above, shifted by 2 rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1047
$wsAll.Range("F5").Value = 280
$wsAll.Range("F6").Value = 2728
$wsAll.Range("F7").Value = 54
